$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"
# Row 4
$ws.Range("E4").Value = 0.425
$ws.Range("F4").Value = 0.07199999999999999
$ws.Range("G4").Value = 0.269
$ws.Range("N4").Value = 0.432
$ws.Range("O4").Value = 0.061
$ws.Range("P4").Value = 0.248
$ws.Range("Q4").Value = 0.024
$ws.Range("R4").Value = 0.017
$ws.Range("S4").Value = 0.131
$ws.Range("W4").Value = 0.295
$ws.Range("X4").Value = 0.11
$ws.Range("Y4").Value = 0.331
$ws.Range("AI4").Value = 0.301
$ws.Range("AJ4").Value = 0.08799999999999999
$ws.Range("AK4").Value = 0.297
$ws.Range("AU4").Value = 0.192
$ws.Range("AV4").Value = 0.029
$ws.Range("AW4").Value = 0.169
$ws.Range("BA4").Value = 2.007
$ws.Range("BB4").Value = 0.157
$ws.Range("BC4").Value = 0.396
$ws.Range("BG4").Value = 0.731
$ws.Range("BH4").Value = 0.139
$ws.Range("BI4").Value = 0.372
$ws.Range("BM4").Value = 0.718
$ws.Range("BN4").Value = 0.075
$ws.Range("BO4").Value = 0.274
$ws.Range("BP4").Value = 0.669
$ws.Range("BQ4").Value = 0.709
# Row 5
$ws.Range("E5").Value = 0.548
$ws.Range("F5").Value = 0.08400000000000001
$ws.Range("G5").Value = 0.29
$ws.Range("N5").Value = 0.748
$ws.Range("O5").Value = 0.076
$ws.Range("P5").Value = 0.276
$ws.Range("Q5").Value = 0.016
$ws.Range("R5").Value = 0.007
$ws.Range("S5").Value = 0.08400000000000001
$ws.Range("W5").Value = 0.285
$ws.Range("X5").Value = 0.11
$ws.Range("Y5").Value = 0.332
$ws.Range("AI5").Value = 0.323
$ws.Range("AJ5").Value = 0.098
$ws.Range("AK5").Value = 0.314
$ws.Range("AU5").Value = 0.375
$ws.Range("AV5").Value = 0.094
$ws.Range("AW5").Value = 0.307
$ws.Range("BA5").Value = 1.357
$ws.Range("BC5").Value = 0.283
$ws.Range("BG5").Value = 0.4
$ws.Range("BH5").Value = 0.05
$ws.Range("BI5").Value = 0.224
$ws.Range("BM5").Value = 0.5610000000000001
$ws.Range("BN5").Value = 0.061
$ws.Range("BO5").Value = 0.247
$ws.Range("BP5").Value = 0.452
$ws.Range("BQ5").Value = 0.461
# Row 6
$ws.Range("E6").Value = 0.479
$ws.Range("N6").Value = 0.548
$ws.Range("Q6").Value = 0.019
$ws.Range("W6").Value = 0.29
$ws.Range("AI6").Value = 0.312
$ws.Range("AU6").Value = 0.254
$ws.Range("BA6").Value = 1.61
$ws.Range("BG6").Value = 0.517
$ws.Range("BM6").Value = 0.63
$ws.Range("BP6").Value = 0.537
$ws.Range("BQ6").Value = 0.555
# Row 7
$ws.Range("E7").Value = 0.518
$ws.Range("N7").Value = 0.653
$ws.Range("Q7").Value = 0.017
$ws.Range("W7").Value = 0.287
$ws.Range("AI7").Value = 0.318
$ws.Range("AU7").Value = 0.315
$ws.Range("BA7").Value = 1.447
$ws.Range("BG7").Value = 0.44
$ws.Range("BM7").Value = 0.587
$ws.Range("BP7").Value = 0.482
$ws.Range("BQ7").Value = 0.494
# Row 8
$ws.Range("E8").Value = 0.606
$ws.Range("F8").Value = 0.112
$ws.Range("G8").Value = 0.334
$ws.Range("N8").Value = 0.781
$ws.Range("O8").Value = 0.062
$ws.Range("P8").Value = 0.248
$ws.Range("Q8").Value = 0.018
$ws.Range("W8").Value = 0.314
$ws.Range("X8").Value = 0.121
$ws.Range("Y8").Value = 0.348
$ws.Range("AI8").Value = 0.345
$ws.Range("AJ8").Value = 0.129
$ws.Range("AK8").Value = 0.36
$ws.Range("AU8").Value = 0.318
$ws.Range("AW8").Value = 0.292
$ws.Range("BA8").Value = 1.756
$ws.Range("BB8").Value = 0.124
$ws.Range("BC8").Value = 0.353
$ws.Range("BG8").Value = 0.5669999999999999
$ws.Range("BH8").Value = 0.106
$ws.Range("BI8").Value = 0.325
$ws.Range("BM8").Value = 0.701
$ws.Range("BN8").Value = 0.062
$ws.Range("BO8").Value = 0.249
$ws.Range("BP8").Value = 0.585
$ws.Range("BQ8").Value = 0.605
# Row 9
$ws.Range("E9").Value = 0.544
$ws.Range("F9").Value = 0.248
$ws.Range("G9").Value = 0.498
$ws.Range("N9").Value = 0.678
$ws.Range("O9").Value = 0.218
$ws.Range("P9").Value = 0.467
$ws.Range("W9").Value = 0.211
$ws.Range("X9").Value = 0.167
$ws.Range("Y9").Value = 0.408
$ws.Range("AI9").Value = 0.267
$ws.Range("AJ9").Value = 0.196
$ws.Range("AK9").Value = 0.442
$ws.Range("BA9").Value = 1.7
$ws.Range("BB9").Value = 0.247
$ws.Range("BC9").Value = 0.497
$ws.Range("BG9").Value = 0.6
$ws.Range("BH9").Value = 0.24
$ws.Range("BI9").Value = 0.49
$ws.Range("BM9").Value = 0.656
$ws.Range("BN9").Value = 0.226
$ws.Range("BO9").Value = 0.475
$ws.Range("BP9").Value = 0.5669999999999999
$ws.Range("BQ9").Value = 0.581
# Row 10
$ws.Range("E10").Value = 0.678
$ws.Range("F10").Value = 0.218
$ws.Range("G10").Value = 0.467
$ws.Range("N10").Value = 0.878
$ws.Range("O10").Value = 0.107
$ws.Range("P10").Value = 0.328
$ws.Range("W10").Value = 0.389
$ws.Range("X10").Value = 0.238
$ws.Range("Y10").Value = 0.487
$ws.Range("AI10").Value = 0.378
$ws.Range("AJ10").Value = 0.235
$ws.Range("AK10").Value = 0.485
$ws.Range("AU10").Value = 0.311
$ws.Range("AV10").Value = 0.214
$ws.Range("AW10").Value = 0.463
$ws.Range("BA10").Value = 2.09
$ws.Range("BB10").Value = 0.244
$ws.Range("BC10").Value = 0.494
$ws.Range("BG10").Value = 0.656
$ws.Range("BH10").Value = 0.226
$ws.Range("BI10").Value = 0.475
$ws.Range("BM10").Value = 0.856
$ws.Range("BN10").Value = 0.124
$ws.Range("BO10").Value = 0.352
$ws.Range("BP10").Value = 0.697
$ws.Range("BQ10").Value = 0.723
# Row 11
$ws.Range("E11").Value = 0.711
$ws.Range("F11").Value = 0.205
$ws.Range("G11").Value = 0.453
$ws.Range("N11").Value = 0.9
$ws.Range("O11").Value = 0.09
$ws.Range("P11").Value = 0.3
$ws.Range("W11").Value = 0.389
$ws.Range("X11").Value = 0.238
$ws.Range("Y11").Value = 0.487
$ws.Range("AI11").Value = 0.411
$ws.Range("AJ11").Value = 0.242
$ws.Range("AK11").Value = 0.492
$ws.Range("AU11").Value = 0.444
$ws.Range("AV11").Value = 0.247
$ws.Range("AW11").Value = 0.497
$ws.Range("BA11").Value = 2.09
$ws.Range("BB11").Value = 0.244
$ws.Range("BC11").Value = 0.494
$ws.Range("BG11").Value = 0.656
$ws.Range("BH11").Value = 0.226
$ws.Range("BI11").Value = 0.475
$ws.Range("BM11").Value = 0.856
$ws.Range("BN11").Value = 0.124
$ws.Range("BO11").Value = 0.352
$ws.Range("BP11").Value = 0.697
$ws.Range("BQ11").Value = 0.726
# Row 12
$ws.Range("E12").Value = 1.422
$ws.Range("F12").Value = 0.775
$ws.Range("G12").Value = 0.88
$ws.Range("N12").Value = 1.482
$ws.Range("O12").Value = 1.069
$ws.Range("P12").Value = 1.034
$ws.Range("W12").Value = 1.629
$ws.Range("X12").Value = 0.576
$ws.Range("Y12").Value = 0.759
$ws.Range("AI12").Value = 1.703
$ws.Range("AJ12").Value = 1.29
$ws.Range("AK12").Value = 1.136
$ws.Range("AU12").Value = 2.714
$ws.Range("AV12").Value = 2.68
$ws.Range("AW12").Value = 1.637
$ws.Range("BA12").Value = 3.728
$ws.Range("BB12").Value = 0.412
$ws.Range("BC12").Value = 0.642
$ws.Range("BG12").Value = 1.102
$ws.Range("BH12").Value = 0.125
$ws.Range("BI12").Value = 0.354
$ws.Range("BM12").Value = 1.299
$ws.Range("BN12").Value = 0.339
$ws.Range("BO12").Value = 0.583
$ws.Range("BP12").Value = 1.243
$ws.Range("BQ12").Value = 1.263
# Row 13
$ws.Range("E13").Value = 1.604
$ws.Range("F13").Value = 0.656
$ws.Range("G13").Value = 0.8100000000000001
$ws.Range("N13").Value = 2.085
$ws.Range("O13").Value = 0.93
$ws.Range("P13").Value = 0.965
$ws.Range("W13").Value = 1.05
$ws.Range("X13").Value = 0.19
$ws.Range("Y13").Value = 0.436
$ws.Range("AI13").Value = 1.28
$ws.Range("AJ13").Value = 0.37
$ws.Range("AK13").Value = 0.608
$ws.Range("AU13").Value = 2.304
$ws.Range("AV13").Value = 0.93
$ws.Range("AW13").Value = 0.964
$ws.Range("BA13").Value = 2.378
$ws.Range("BB13").Value = 0.302
$ws.Range("BC13").Value = 0.55
$ws.Range("BG13").Value = 0.59
$ws.Range("BH13").Value = 0.07199999999999999
$ws.Range("BI13").Value = 0.268
$ws.Range("BM13").Value = 0.908
$ws.Range("BN13").Value = 0.284
$ws.Range("BO13").Value = 0.533
$ws.Range("BP13").Value = 0.793
$ws.Range("BQ13").Value = 0.733
